$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need an explicit Text
# number format first, otherwise Excel auto-converts the assigned string
# into a numeric value instead of keeping it as text (the source data is
# all text/inlineStr cells).
$textGuardCells = @("D5","D6","D14","D19","D21","D22","D23","D24","D25","D28","D29","D30","D31","D32","D33","D34","D36","D37","D38","D40","D43","D45","D46","D48","D49","D50")
foreach ($cellRef in $textGuardCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "63.145.76"
$ws.Range("E2").Value = "  -1.07%  "
$ws.Range("D3").Value = "2.627.86"
$ws.Range("E3").Value = "  -1.26%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "601.13"
$ws.Range("E5").Value = "  +1.05%  "
$ws.Range("D6").Value = "146.20"
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  -0.88%  "
$ws.Range("D9").Value = "2.628.02"
$ws.Range("E9").Value = "  -1.21%  "
$ws.Range("E10").Value = "  -0.47%  "
$ws.Range("E11").Value = "  -0.97%  "
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("E13").Value = "  +1.81%  "
$ws.Range("D14").Value = "27.13"
$ws.Range("E14").Value = "  -1.72%  "
$ws.Range("D15").Value = "3.099.71"
$ws.Range("E15").Value = "  -1.08%  "
$ws.Range("D16").Value = "63.015.38"
$ws.Range("E17").Value = "  -2.02%  "
$ws.Range("D18").Value = "2.616.07"
$ws.Range("E18").Value = "  -1.60%  "
$ws.Range("D19").Value = "11.29"
$ws.Range("E19").Value = "  -1.17%  "
$ws.Range("E20").Value = "  +2.40%  "
$ws.Range("D21").Value = "339.89"
$ws.Range("E21").Value = "  -1.10%  "
$ws.Range("D22").Value = "6.87"
$ws.Range("E22").Value = "  +1.09%  "
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("D24").Value = "5.57"
$ws.Range("E24").Value = "  -3.88%  "
$ws.Range("D25").Value = "66.46"
$ws.Range("E25").Value = "  -2.41%  "
$ws.Range("E26").Value = "  -3.32%  "
$ws.Range("E27").Value = "  +0.60%  "
$ws.Range("B28").Value = "Bittensor"
$ws.Range("C28").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D28").Value = "548.77"
$ws.Range("E28").Value = "  -4.60%  "
$ws.Range("B29").Value = "SuiNetwork"
$ws.Range("C29").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D29").Value = "1.51"
$ws.Range("E29").Value = "  -7.68%  "
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").Value = "0.163"
$ws.Range("E30").Value = "  -2.48%  "
$ws.Range("D31").Value = "0.995"
$ws.Range("E31").Value = "  -0.64%  "
$ws.Range("D32").Value = "7.84"
$ws.Range("E32").Value = "  -2.98%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").Value = "2.00"
$ws.Range("E33").Value = "  -0.27%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "1.74"
$ws.Range("E34").Value = "  -2.89%  "
$ws.Range("D35").Value = "0.0₃0801"
$ws.Range("E35").Value = "  -2.28%  "
$ws.Range("D36").Value = "5.29"
$ws.Range("E36").Value = "  +11.76%  "
$ws.Range("D37").Value = "165.83"
$ws.Range("E37").Value = "  -5.40%  "
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("E39").Value = "  +0.15%  "
$ws.Range("D40").Value = "18.94"
$ws.Range("E40").Value = "  -1.27%  "
$ws.Range("E41").Value = "  +6.01%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("D43").Value = "167.48"
$ws.Range("E43").Value = "  -2.52%  "
$ws.Range("E44").Value = "  -1.08%  "
$ws.Range("D45").Value = "22.31"
$ws.Range("E45").Value = "  +2.17%  "
$ws.Range("D46").Value = "0.0567"
$ws.Range("E46").Value = "  +2.41%  "
$ws.Range("E47").Value = "  -1.55%  "
$ws.Range("D48").Value = "0.0242"
$ws.Range("E48").Value = "  +0.82%  "
$ws.Range("D49").Value = "0.0957"
$ws.Range("E49").Value = "  -0.49%  "
$ws.Range("D50").Value = "18.56"
$ws.Range("E50").Value = "  -0.93%  "
$ws.Range("E51").Value = "  +0.12%  "
